{"js": "// Add \"Hello World\" text to the last (empty) paragraph of the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"Hello World\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add \"Hello World\" text to the last (empty) paragraph of the document body.\n$d = $word.ActiveDocument\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.Text = \"Hello World\"\n"}
